# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This script updates the DAMSLTag (column I) and DialogAct (column J) values
# for the rows whose automatic dialog-act annotation changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 19;  I = "ba"; J = "Appreciation" },
    @{ Row = 20;  I = "ba"; J = "Appreciation" },
    @{ Row = 28;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 37;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 75;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 94;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 101; I = "sv"; J = "Statement-opinion" },
    @{ Row = 118; I = "aa"; J = "Agree/Accept" },
    @{ Row = 134; I = "aa"; J = "Agree/Accept" },
    @{ Row = 139; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 140; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 150; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 155; I = "aa"; J = "Agree/Accept" },
    @{ Row = 158; I = "aa"; J = "Agree/Accept" },
    @{ Row = 173; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 182; I = "%";  J = "Uninterpretable" },
    @{ Row = 186; I = "aa"; J = "Agree/Accept" },
    @{ Row = 187; I = "sv"; J = "Statement-opinion" },
    @{ Row = 192; I = "sv"; J = "Statement-opinion" },
    @{ Row = 225; I = "sv"; J = "Statement-opinion" },
    @{ Row = 227; I = "sv"; J = "Statement-opinion" },
    @{ Row = 234; I = "%";  J = "Uninterpretable" },
    @{ Row = 238; I = "sv"; J = "Statement-opinion" },
    @{ Row = 247; I = "sv"; J = "Statement-opinion" },
    @{ Row = 261; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 267; I = "sv"; J = "Statement-opinion" },
    @{ Row = 275; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 289; I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 294; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 296; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 306; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 308; I = "aa"; J = "Agree/Accept" },
    @{ Row = 326; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 330; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 352; I = "sv"; J = "Statement-opinion" },
    @{ Row = 356; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 364; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 371; I = "sv"; J = "Statement-opinion" },
    @{ Row = 375; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
